$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BTS NTS Modal Profile Data")

# Insert a new row at position 37 (pushes existing row 37 "AVLo-freight" header and below down to 38+)
$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the weighted-value-adjusted label and formula
$ws.Range("A37").Value = "weighted value, adjusted for number of train cars per locomotive"
$ws.Range("B37").Formula = "=B36/10"

# Match style of row 36 (A37 uses style similar to A36's "s=10" -- label style; B37 uses "s=6" like B25 etc.)
$ws.Range("A37").WrapText = $true
$ws.Range("B37").NumberFormat = "0"
$ws.Rows.Item(37).RowHeight = 16
$ws.Rows.Item(36).RowHeight = 16
$ws.Rows.Item(60).RowHeight = 16

# Row 36's formula is unaffected by the insertion (it's above it), but the "AVLo-passengers"
# sheet should now reference the new adjusted row 37 instead of the raw row 36 value.
$wsPax = $wb.Worksheets.Item("AVLo-passengers")
$wsPax.Range("B5").Formula = "='BTS NTS Modal Profile Data'!B37"

# Update view/selection state to match the saved workbook session:
#  - "About" is no longer the active tab
#  - "BTS NTS Modal Profile Data" view is scrolled down with C34 selected
#  - "AVLo-passengers" becomes the active tab, with F12 selected
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A44").Select() | Out-Null

$ws.Activate()
$ws.Range("C34").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 14

$wsPax.Activate()
$wsPax.Range("F12").Select() | Out-Null

# Header row heights adjusted on the two "AVLo" sheets
$wsPax.Rows.Item(1).RowHeight = 48

$wsFreight = $wb.Worksheets.Item("AVLo-freight")
$wsFreight.Rows.Item(1).RowHeight = 32

$wb.Save()
